# Optimize security vulnerability checks
# Appends one new trailing data row to each of the four log sheets
# (ROW50-FE-LIFTER, ROW50-MID-LIFTER, ROW11-FE-LIFTER, ROW11-MID-LIFTER),
# mirroring the existing "time / len / ID / actual-len / checksum (+ _DEC)"
# record layout already present in each sheet.

$wb = $excel.ActiveWorkbook

$dateFormat = "YYYY-MM-DD HH:MM:SS"

function Add-LogRow {
    param($ws, $row, $timeValue, $col2, $col3, $col4, $col5, $col6, $col7, $col7IsText, $col8, $col9)

    # Column A: timestamp, keep the same date/time number format as the rest
    # of the column (re-using cellXf index 2 rather than minting a new one).
    $ws.Cells.Item($row, 1).NumberFormat = $dateFormat
    $ws.Cells.Item($row, 1).Value = $timeValue

    # Columns B-E: raw hex byte strings (stored as text).
    $ws.Cells.Item($row, 2).Value = $col2
    $ws.Cells.Item($row, 3).Value = $col3
    $ws.Cells.Item($row, 4).Value = $col4
    $ws.Cells.Item($row, 5).Value = $col5

    # Column F: decimal length.
    $ws.Cells.Item($row, 6).Value = $col6

    # Column G: decimal ID. Usually numeric, but one sheet keeps it as exact
    # text since the true integer has more digits than a double can hold.
    if ($col7IsText) {
        $ws.Cells.Item($row, 7).NumberFormat = "@"
        $ws.Cells.Item($row, 7).Value = $col7
        $ws.Cells.Item($row, 7).Style = "Normal"
    } else {
        $ws.Cells.Item($row, 7).Value = $col7
    }

    # Columns H-I: decimal actual-length / checksum.
    $ws.Cells.Item($row, 8).Value = $col8
    $ws.Cells.Item($row, 9).Value = $col9
}

# --- ROW50-FE-LIFTER: new row 77 ---
$ws1 = $wb.Worksheets.Item("ROW50-FE-LIFTER")
Add-LogRow $ws1 77 45762.76107120371 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," `
    "0x01,0x4a" `
    "0xe" `
    400 `
    568631262647114000000000.0 $false `
    330 `
    14

# --- ROW50-MID-LIFTER: new row 79 ---
$ws2 = $wb.Worksheets.Item("ROW50-MID-LIFTER")
Add-LogRow $ws2 79 45762.72451388889 `
    "0x01,0x90 " `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x20,0x41,0x0c," `
    "0x01,0x52" `
    "0x19" `
    400 `
    "568631262647113771663628" $true `
    338 `
    25

# --- ROW11-FE-LIFTER: new row 77 ---
$ws3 = $wb.Worksheets.Item("ROW11-FE-LIFTER")
Add-LogRow $ws3 77 45762.79501923611 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x1a,0x41,0x0c," `
    "0x01,0x4a" `
    "0x14" `
    400 `
    568631262647114000000000.0 $false `
    330 `
    20

# --- ROW11-MID-LIFTER: new row 77 ---
$ws4 = $wb.Worksheets.Item("ROW11-MID-LIFTER")
Add-LogRow $ws4 77 45762.92137438658 `
    "0x01,0x90" `
    "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x11,0x40,0x0c," `
    "0x01,0x52" `
    "0x19" `
    400 `
    568631262647114000000000.0 $false `
    338 `
    25
